# Insert a new weekly record for "Terminal La Palmera de La Serena - Zanahoria".
# This shifts the existing rows 426..456 down to 427..457 and inserts a brand
# new row 426 with the latest week's data (dimension grows from R456 to R457).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push everything from row 426 down by one row, carrying the row's formatting
# (this keeps the date-style on column D for the new blank row).
$ws.Rows.Item(426).Insert()

# Populate the newly inserted row 426 with the new weekly observation.
$ws.Cells.Item(426, 1).Value = 8
$ws.Cells.Item(426, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(426, 3).Value = "Coquimbo"
$ws.Cells.Item(426, 4).Value = 44931
$ws.Cells.Item(426, 5).Value = 4
$ws.Cells.Item(426, 6).Value = 100114013
$ws.Cells.Item(426, 7).Value = "Zanahoria"
$ws.Cells.Item(426, 8).Value = "Sin especificar"
$ws.Cells.Item(426, 9).Value = "Primera"
$ws.Cells.Item(426, 10).Value = 600
$ws.Cells.Item(426, 11).Value = 5800
$ws.Cells.Item(426, 12).Value = 6000
$ws.Cells.Item(426, 13).Value = 5900
$ws.Cells.Item(426, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(426, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(426, 16).Value = 295
$ws.Cells.Item(426, 17).Value = 20
$ws.Cells.Item(426, 18).Value = "Hortaliza"
